$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A width (target raw width 15.42578125; engine quantizes
# ColumnWidth to the nearest 1/6, so feed the value that rounds closest)
$ws.Columns.Item(1).ColumnWidth = 14.666666666666666

# Update values in A1:A33
$values = @(
    0.20264560746247184,
    -0.057644320437837848,
    -0.0039999999827529109,
    -0.0079999999682023315,
    -0.002999999983220647,
    -0.0019999999834272586,
    -0.0099999999572650822,
    0.0026423739155467807,
    -0.0019999999829685144,
    -0.0019999999827344794,
    -0.0029999999794512178,
    -0.0034999999778393409,
    -0.0034999999781257785,
    -0.0079999999636273245,
    0.020193375024809335,
    -0.0019999999834552362,
    -0.0019999999831137316,
    -0.0039999999764495087,
    -0.003999999985934366,
    -0.0039999999849467116,
    -0.0039999999847735168,
    -0.0039999999846420664,
    -0.0049999999783354454,
    -0.019999999927158285,
    -0.019999999926148426,
    -0.002499999981363743,
    -0.0024999999807668871,
    -0.0019999999798177015,
    -0.0069999999615575348,
    -0.059999999787059277,
    0.042229998377326439,
    -0.0099999999522566441,
    -0.0039999999718425272
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
